# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.374.80'
$ws.Range("E2").Value = '  +3.32%  '

$ws.Range("D3").Value = '1.863.92'
$ws.Range("E3").Value = '  +1.97%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.13'
$ws.Range("E5").Value = '  +1.91%  '

$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4706'
$ws.Range("E7").Value = '  +2.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3969'
$ws.Range("E8").Value = '  +3.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.56'
$ws.Range("E9").Value = '  +2.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08010'
$ws.Range("E10").Value = '  +1.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9958'
$ws.Range("E11").Value = '  +2.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.96'
$ws.Range("E12").Value = '  +4.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.024'
$ws.Range("E13").Value = '  +2.46%  '

$ws.Range("D14").Value = '1.864.87'
$ws.Range("E14").Value = '  +1.93%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.240'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.37'
$ws.Range("E16").Value = '  +2.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001039'
$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06630'
$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.51'
$ws.Range("E20").Value = '  +1.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  -0.36%  '

$ws.Range("D22").Value = '28.387.88'
$ws.Range("E22").Value = '  +3.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.459'
$ws.Range("E23").Value = '  +2.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  +2.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.270'
$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").Value = '2.081.99'
$ws.Range("E26").Value = '  +1.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.93'
$ws.Range("E27").Value = '  +2.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.70'
$ws.Range("E28").Value = '  +1.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.111'
$ws.Range("E29").Value = '  +2.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.459'
$ws.Range("E30").Value = '  +4.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.48'
$ws.Range("E31").Value = '  +1.10%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09501'
$ws.Range("E32").Value = '  +2.31%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9596'
$ws.Range("E33").Value = '  +1.01%  '

$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.375'
$ws.Range("E35").Value = '  +4.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.348'
$ws.Range("E36").Value = '  +2.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06108'
$ws.Range("E37").Value = '  +2.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02245'
$ws.Range("E38").Value = '  +2.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.271'
$ws.Range("E39").Value = '  +3.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.176'
$ws.Range("E40").Value = '  +1.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5916'
$ws.Range("E41").Value = '  +2.24%  '

$ws.Range("E42").Value = '  -0.44%  '

$ws.Range("E43").Value = '  +1.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.28'
$ws.Range("E44").Value = '  +2.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.258'
$ws.Range("E45").Value = '  -1.49%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07451'
$ws.Range("E46").Value = '  +12.11%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5535'
$ws.Range("E47").Value = '  +0.83%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.10'
$ws.Range("E48").Value = '  +1.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.945'
$ws.Range("E49").Value = '  +4.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.062'
$ws.Range("E50").Value = '  +13.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.84'
$ws.Range("E51").Value = '  +1.78%  '
